$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format before writing numeric-looking strings,
# so Excel does not auto-convert them to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.710.71'
$ws.Range('E2').Value = '  +1.08%  '

$ws.Range('D3').Value = '1.834.90'
$ws.Range('E3').Value = '  +1.62%  '

$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.30%  '

$ws.Range('D5').Value = '309.06'
$ws.Range('E5').Value = '  +0.70%  '

$ws.Range('E6').Value = '  +0.15%  '

$ws.Range('D7').Value = '0.4671'
$ws.Range('E7').Value = '  +3.35%  '

$ws.Range('D8').Value = '0.3611'
$ws.Range('E8').Value = '  +0.47%  '

$ws.Range('D9').Value = '0.07153'

$ws.Range('D10').Value = '0.9327'
$ws.Range('E10').Value = '  +4.79%  '

$ws.Range('D11').Value = '19.51'
$ws.Range('E11').Value = '  +0.21%  '

$ws.Range('D12').Value = '0.07667'
$ws.Range('E12').Value = '  -1.79%  '

$ws.Range('D13').Value = '1.858.47'
$ws.Range('E13').Value = '  +2.71%  '

$ws.Range('D14').Value = '5.262'
$ws.Range('E14').Value = '  -0.55%  '

$ws.Range('D15').Value = '6.363'
$ws.Range('E15').Value = '  +0.52%  '

$ws.Range('D16').Value = '87.79'
$ws.Range('E16').Value = '  +3.19%  '

$ws.Range('E17').Value = '  +0.29%  '

$ws.Range('D18').Value = '0.000008563'
$ws.Range('E18').Value = '  +1.05%  '

$ws.Range('D19').Value = '1.006'
$ws.Range('E19').Value = '  +0.15%  '

$ws.Range('D20').Value = '26.729.78'
$ws.Range('E20').Value = '  +0.98%  '

$ws.Range('D21').Value = '14.28'
$ws.Range('E21').Value = '  +0.03%  '

$ws.Range('D22').Value = '5.020'
$ws.Range('E22').Value = '  +0.94%  '

$ws.Range('E23').Value = '  +0.67%  '

$ws.Range('D24').Value = '1.913'
$ws.Range('E24').Value = '  -2.51%  '

$ws.Range('D25').Value = '151.86'
$ws.Range('E25').Value = '  +0.73%  '

$ws.Range('D26').Value = '17.97'
$ws.Range('E26').Value = '  +0.76%  '

$ws.Range('D27').Value = '2.005'
$ws.Range('E27').Value = '  -2.30%  '

$ws.Range('D28').Value = '113.75'
$ws.Range('E28').Value = '  +1.62%  '

$ws.Range('D29').Value = '4.894'
$ws.Range('E29').Value = '  +0.64%  '

$ws.Range('D30').Value = '0.08823'
$ws.Range('E30').Value = '  +1.47%  '

$ws.Range('D31').Value = '3.153'
$ws.Range('E31').Value = '  +1.29%  '

$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.182'
$ws.Range('E32').Value = '  +6.90%  '

$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '2.820'
$ws.Range('E33').Value = '  -0.88%  '

$ws.Range('D34').Value = '0.7413'
$ws.Range('E34').Value = '  +2.93%  '

$ws.Range('D35').Value = '4.450'
$ws.Range('E35').Value = '  +0.16%  '

$ws.Range('E36').Value = '  +1.03%  '

$ws.Range('D37').Value = '2.970'
$ws.Range('E37').Value = '  +2.38%  '

$ws.Range('E38').Value = '  -0.82%  '

$ws.Range('D39').Value = '0.05149'
$ws.Range('E39').Value = '  +1.03%  '

$ws.Range('D40').Value = '6.912'
$ws.Range('E40').Value = '  +1.94%  '

$ws.Range('D41').Value = '0.5074'
$ws.Range('E41').Value = '  -1.01%  '

$ws.Range('E42').Value = '  -0.70%  '

$ws.Range('D43').Value = '8.130'

$ws.Range('D44').Value = '0.4671'
$ws.Range('E44').Value = '  +0.42%  '

$ws.Range('D45').Value = '1.006'
$ws.Range('E45').Value = '  +0.15%  '

$ws.Range('D46').Value = '10.16'
$ws.Range('E46').Value = '  +1.70%  '

$ws.Range('D47').Value = '99.18'
$ws.Range('E47').Value = '  -1.19%  '

$ws.Range('D48').Value = '1.577'
$ws.Range('E48').Value = '  +0.34%  '

$ws.Range('D49').Value = '0.06033'

$ws.Range('D50').Value = '64.08'
$ws.Range('E50').Value = '  +0.38%  '

$ws.Range('D51').Value = '36.02'
$ws.Range('E51').Value = '  -0.25%  '

# Restore default cell style on column D (remove the temporary text format)
$dRange.Style = "Normal"
